# Apply the "Added exclude user list, irrelevant filter, picture display feature"
# change to the comments.xlsx workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("mainComments")
$ws2 = $wb.Worksheets.Item("replys")

# --- New "irrelevantTag" / "image" columns (headers) ---
$ws1.Range("G1").Value = "irrelevantTag"
$ws1.Range("H1").Value = "image"
$ws2.Range("G1").Value = "irrelevantTag"
$ws2.Range("H1").Value = "image"

# --- Flag existing comments/replies that are irrelevant to the topic ---
$ws1.Range("G3").Value = 1
$ws1.Range("G4").Value = 1

$ws2.Range("G3").Value = 1
$ws2.Range("G5").Value = 1

# --- New comment row (mainComments!7): an off-topic reply from a new user ---
$ws1.Range("A7").Value = 6
$ws1.Range("B7").Value = "阿扁"
$ws1.Range("C7").Value = "三小"
$ws1.Range("D7").Value = "2025-05-27 10:16"
$ws1.Range("E7").Value = "🧑‍🦱"
$ws1.Range("F7").Value = 3

# --- New reply row (replys!6): same new user replying in the replys sheet ---
$ws2.Range("A6").Value = 6
$ws2.Range("B6").Value = "阿扁"
$ws2.Range("C6").Value = "蛤"
$ws2.Range("D6").Value = "2025-05-27 10:16"
$ws2.Range("E6").Value = "🧑‍🦱"
$ws2.Range("F6").Value = 3

# Give the new commenter's name cells the new "Microsoft JhengHei" font style
$ws1.Range("B7").Font.Name = "Microsoft JhengHei"
$ws1.Range("B7").VerticalAlignment = -4108
$ws2.Range("B6").Font.Name = "Microsoft JhengHei"
$ws2.Range("B6").VerticalAlignment = -4108

# --- Picture-display feature: attach image filenames to comments/replies ---
$ws1.Range("H7").Value = "pic3.jpg"
$ws2.Range("H2").Value = "pic2.jpg"
$ws2.Range("H6").Value = "pic4.jpg"

# --- Selection / active sheet bookkeeping to mirror the authored workbook ---
$ws1.Range("H7").Select()
$ws2.Range("H6").Select()
$ws2.Activate()
